$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Last minute cancellation" - remove Minna Ahokas's registration row (row 52)
# and shift everything below it up by one row.
$ws.Rows.Item(52).Delete()

# Column A holds a literal running index ("#"); deleting a row does not
# renumber these hard-coded values automatically, so fix them up for all
# the rows that shifted (former #52.."#69" entries, now sitting in rows
# 52-69, need their number reduced by one: 51..68).
$lastRow = 69
for ($r = 52; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# Restore the view state captured the next time the workbook was saved.
$wb.Windows.Item(1).WindowState = $wb.Windows.Item(1).WindowState
$excel.ActiveWindow.Left = 7440
$excel.ActiveWindow.Top = 2020

$ws.Range("A49:A69").Select()
$excel.ActiveWindow.ScrollRow = 60
